# Update the two "read the observation..." GitHub URL cells with the
# actual (now de-anonymised) repository links.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "https://github.com/diegocedrim/argouml"
$ws.Range("B6").Value = "https://github.com/diegocedrim/xerces"

# These two cells used to carry a bold, wrapped "observation" style;
# restyle them to match the rest of the GitHub URL column (regular
# weight, vertically centered, wrapped).
$ws.Range("B5:B6").Font.Bold = $false
$ws.Range("B5:B6").VerticalAlignment = -4108
$ws.Range("B5:B6").WrapText = $true

# Remove the footnote row: unmerge A28:F28 and clear its text.
$ws.Range("A28:F28").UnMerge()
$ws.Range("A28").ClearContents()

# Move the active selection to where the author left off editing.
$ws.Range("D22").Select()
